$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build each row as its own array, then combine into a jagged array
$r1 = @("Names", "Tuesday", "Wednesday", "Thursday", "Friday", "Saturday", "Sunday", "Monday")
$r2 = @("Pool Hours", "10:45-8", "10:45-8", "10:45-8", "10:45-9", "10:30-8", "11:00-8", "closed")
$r3 = @("", "July2nd", "Jully3rd", "July4th", "July5th", "July6th", "July7th", "July8th")
$r4 = @("Barry Ray", "X", "X", "X", "X", "X", "X", "X")
$r5 = @("Blake Butz", "3:30-8", "10:15-3:30", "10:15-3:30", "X", "X", "X", "X")
$r6 = @("Kate North", "10:15-3:30", "10:15-3:30", "10:00-3", "X", "3:30-8", "3:30-8", "X")
$r7 = @("Emerson Metzger", "10:15-3:30", "OFF", "3:00-7", "X", "X", "X", "X")
$r8 = @("Avery Larsen", "3:30-8", "3:30-8", "10:00-3", "10:15-4", "3:30-8", "3:30-8", "X")
$r9 = @("Austin Page", "OFF", "OFF", "X", "4:00-9", "10:00-3:30", "3:30-8", "X")
$r10 = @("Riley White ", "3:30-8", "10:15-3:30", "10:00-3", "10:15-4", "X", "3:30-8", "X")
$r11 = @("Robert Wade", "10:15-3:30", "10:30-3:30", "X", "X", "3:30-8", "10:30-3:30", "X")
$r12 = @("Tatum Plunk", "10:30-3:30", "X", "X", "X", "X", "X", "X")
$r13 = @("Michael Vangruber", "3:30-8", "X", "X", "X", "10:00-3:30", "10:30-3:30", "X")
$r14 = @("Jackson Blakely", "10:30-3:30", "OFF", "3:00-7", "4:00-9", "1:00-6:00", "10:30-3:30", "X")
$r15 = @("Addison Clark", "10:30-3:30", "10:30-3:30", "10:00-3", "10:15-4", "3:30-8", "1:00-6:00", "X")
$r16 = @("Madison Johnson", "OFF", "3:30-8", "3:00-7", "10:30-4", "10:00-3:30", "10:45-3:30", "X")
$r17 = @("Nathan Debergh", "X", "X", "X", "X", "X", "X", "X")
$r18 = @("Asher Bobbett", "OFF", "X", "X", "X", "X", "X", "X")
$r19 = @("Blake Ucherek", "X", "X", "X", "X", "X", "3:30-8", "X")
$r20 = @("Ethan Van Horn ", "OFF", "3:30-8", "3:00-7", "X", "X", "X", "X")
$r21 = @("Kai King", "OFF", "3:30-8", "11:00-5", "4:00-9", "3:30-8", "OFF", "X")
$r22 = @("Madeline Ellison", "3:30-8", "X", "X", "4:00-9", "OFF", "10:45-3:30", "X")
$r23 = @("Tyler Carpenter", "OFF", "OFF", "3:00-7", "4:00-9", "3:30-8", "OFF", "X")
$r24 = @("Holden ", "X", "X", "X", "X", "X", "X", "X")
$r25 = @("Jayden Garcia ", "3:30-8", "3:30-8", "3:00-7", "10:30-4", "10:15-3:30", "3:30-8", "X")
$r26 = @("Naya Okonkwo", "OFF", "10:30-3:30", "3:00-7", "4:00-9", "10:15-3:30", "10:45-3:30", "X")
$r27 = @("Bella Hamilton", "X", "X", "12:00-6", "10:30-4", "X", "OFF", "X")
$r28 = @("Phillip Thompson", "OFF", "3:30-8", "10:00-3", "1:00-6:00", "10:15-3:30", "OFF", "X")
$r29 = @("Brent Horwitz", "X", "X", "X", "X", "X", "OFF", "X")

$data = @($r1, $r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16, $r17, $r18, $r19, $r20, $r21, $r22, $r23, $r24, $r25, $r26, $r27, $r28, $r29)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
